$wb = $excel.ActiveWorkbook

# Row 5 on every sheet corresponds to file "9f897edb-7a44-463d-be91-d41781ed7fbc.md"
# whose Status moves from "Ready for handoff" to "In Translation".

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B5").Value = "In Translation"
$overview.Range("C5").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B5").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B5").Value = "In Translation"
